# docs(Slide): Submissão do Documento
#
# 1) Slide 1 - "Rounded Rectangle 2": merge the two runs of the second
#    paragraph ("Guilherme Gois " + "Cruz Coelho RA: 25.00702-2") into a
#    single run, keeping the formatting (dirty="0") of the second run.
# 2) Slide 2 - last picture ("Imagem 27"): tiny resize/reposition.
# 3) Slide 8 - "Text 3": split "O teste incluiria um representante..."
#    so "incluiria" becomes "incluiu" as its own run.
# 4) Slide 8 - "Text 4": split "Os usuários realizariam tarefas..." so
#    "usuários" is isolated and "realizariam" becomes "realizaram".
# 5) Slide 8 - "Text 5": split "Seriam coletadas métricas..." so
#    "Seriam" becomes "Foram" as its own run.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 1: merge "Guilherme Gois " + "Cruz Coelho RA: 25.00702-2"
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$nameShape = $s1.Shapes.Item(2)
$nameTr = $nameShape.TextFrame.TextRange

# Paragraph 2 currently spans chars 39-79 ("Guilherme Gois " + "Cruz Coelho RA: 25.00702-2").
# Insert the "Guilherme Gois " prefix into the second run (which carries dirty="0"),
# then blank out the original first-run text so only one run remains.
$secondRun = $nameTr.Characters(54, 26)
[void]$secondRun.InsertBefore("Guilherme Gois ")
$firstRun = $nameTr.Characters(39, 15)
$firstRun.Text = ""

# ---------------------------------------------------------------------
# 2) Slide 2: nudge the full-bleed picture "Imagem 27"
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$bannerPic = $s2.Shapes.Item(28)
$bannerPic.Left = -0.51386
$bannerPic.Width = 1152.639

# ---------------------------------------------------------------------
# 3) Slide 8: "Text 3" -> "incluiria" becomes "incluiu"
# ---------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$publico = $s8.Shapes.Item(4).TextFrame.TextRange
$publico.Characters(19, 9).Text = "incluiu"

# ---------------------------------------------------------------------
# 4) Slide 8: "Text 4" -> isolate "usuários" and "realizariam" -> "realizaram"
# ---------------------------------------------------------------------
$tarefas = $s8.Shapes.Item(5).TextFrame.TextRange
$tarefas.Characters(14, 8).Text = "usuários"
$tarefas.Characters(23, 11).Text = "realizaram"

# ---------------------------------------------------------------------
# 5) Slide 8: "Text 5" -> "Seriam" becomes "Foram"
# ---------------------------------------------------------------------
$metricas = $s8.Shapes.Item(6).TextFrame.TextRange
$metricas.Characters(12, 6).Text = "Foram"
